$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the header formatting
# already used by the other header cells (bold, bordered, centered).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for the new columns I (I0) and J (IF), rows 2-23.
$data = @{
    2  = @(8, 8)
    3  = @(9, 9)
    4  = @(9, 9)
    5  = @(7, 7)
    6  = @(8, 8)
    7  = @(6, 6)
    8  = @(6, 7)
    9  = @(7, 7)
    10 = @(4, 5)
    11 = @(9, 9)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(7, 7)
    15 = @(6, 6)
    16 = @(7, 7)
    17 = @(8, 8)
    18 = @(8, 8)
    19 = @(7, 7)
    20 = @(7, 7)
    21 = @(6, 6)
    22 = @(6, 6)
    23 = @(9, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
